$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.455.76"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.104.99"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.18"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4615"
$ws.Range("E8").Value = "  +6.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.35"
$ws.Range("E9").Value = "  +12.44%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.44"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.095.43"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.800"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.947"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.47"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06631"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.29"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.291"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.509.08"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.362"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.335.84"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.32"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.565"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.92"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.78"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.694"
$ws.Range("E33").Value = "  +9.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.164"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.926"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("E36").Value = "  +7.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02574"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06830"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.557"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.84"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2291"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6895"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.340"
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6391"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.96"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.666"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +25.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.248"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.44"
$ws.Range("E51").Value = "  +0.92%  "
